# Apply the committed change to TestCases.xlsx:
#  - Update sheet1 ("leathershoplogin"): the last step now does a
#    verifyExactText (was verifyText) and its expected value moves from
#    column E into a new column F ("Verification Text").
#  - Add a second worksheet ("cangotomensformalshoepage") with a new keyword
#    driven test case (TC_02), make it the active sheet/tab.
#
# NOTE: cell values are written in the same order the original author must
# have typed them (inferred from the shared-string table ordering) so that
# new shared strings line up the same way.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: "leathershoplogin" -------------------------------------------

# Row 6 ("Verify log in") used to do a plain verifyText with the expected
# value in column E. It now does a verifyExactText.
$ws1.Range("D6").Value = "verifyExactText"

# Add the new "Verification Text" header in column F.
$ws1.Range("F1").Value = "Verification Text"

# Move the expected value from E6 to F6.
$expected = $ws1.Range("E6").Value()
$ws1.Range("F6").Value = $expected
$ws1.Range("E6").ClearContents()

# The active cell/tab moves to the new sheet; sheet1 keeps a header-row
# selection.
$ws1.Rows.Item(1).Select()

# --- Sheet2: "cangotomensformalshoepage" -----------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "cangotomensformalshoepage"

# Header row (same headers as sheet1, including the new "Verification Text").
$ws2.Range("A1").Value = "TC_NO"
$ws2.Range("B1").Value = "Steps"
$ws2.Range("C1").Value = "Locator Key"
$ws2.Range("D1").Value = "Actions"
$ws2.Range("E1").Value = "Test Data"
$ws2.Range("F1").Value = "Verification Text"
$ws2.Range("A1:F1").Font.Bold = $true

# TC_02 test-case rows (column D of row 2 is filled in last, below).
$ws2.Range("A2").Value = "TC_02"
$ws2.Range("B2").Value = "Log in to site"

$ws2.Range("A3").Value = "TC_02"
$ws2.Range("B3").Value = "Move to menu option men"
$ws2.Range("C3").Value = "leathershop.menu.men"
$ws2.Range("D3").Value = "moveToElement"

$ws2.Range("A4").Value = "TC_02"
$ws2.Range("B4").Value = "Move to sub menu formal and click"
$ws2.Range("C4").Value = "leathershop.men.submenu.formal"
$ws2.Range("D4").Value = "moveToElementAndClick"

$ws2.Range("A5").Value = "TC_02"
$ws2.Range("B5").Value = "Vefify formal page"
$ws2.Range("C5").Value = "leathershop.men.formalshoepage.formaltext"
$ws2.Range("D5").Value = "verifyTextContains"
$ws2.Range("F5").Value = "FORMAL"

$ws2.Range("A6").Value = "TC_02"
$ws2.Range("B6").Value = "Log out from site"
$ws2.Range("D6").Value = "logoutFromSite"

$ws2.Range("D2").Value = "loginToSite"

# View state: select C12 on the new sheet (matches target selection) and
# make it the active/visible tab.
$ws2.Range("C12").Select()
$ws2.PageSetup.Orientation = 1
$ws2.Activate()
